$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2603
$ws.Range("I15").Value = 2603
$ws.Range("K15").Value = 7809
$ws.Range("M15").Value = -7640

$ws.Range("H17").Value = 554.53
$ws.Range("J17").Value = 547.6421
$ws.Range("L17").Value = 1642.9263
$ws.Range("N17").Value = -1978.9263

$ws.Range("H32").Value = 3621.25
$ws.Range("I32").Value = 2499
$ws.Range("J32").Value = 3995.3333
$ws.Range("K32").Value = 2499
$ws.Range("L32").Value = 3995.3333
$ws.Range("M32").Value = -2173
$ws.Range("N32").Value = -4647.3333

$ws.Range("H39").Value = 1394.1
$ws.Range("I39").Value = 407.33334
$ws.Range("J39").Value = 2874.25
$ws.Range("K39").Value = 1222.00002
$ws.Range("L39").Value = 8622.75
$ws.Range("M39").Value = -926.0000199999999
$ws.Range("N39").Value = -9214.75

$ws.Range("H98").Value = 1514.091
$ws.Range("I98").Value = 1415.5
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 1415.5
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = 82.5
$ws.Range("N98").Value = -5496

$ws.Range("H101").Value = 334.4
$ws.Range("I101").Value = 343
$ws.Range("J101").Value = 300
$ws.Range("K101").Value = 1029
$ws.Range("L101").Value = 900
$ws.Range("M101").Value = 593
$ws.Range("N101").Value = -4144

$ws.Range("H112").Value = 1029.1666
$ws.Range("J112").Value = 1054.7812
$ws.Range("L112").Value = 3164.3436
$ws.Range("N112").Value = -5380.3436

$ws.Range("H122").Value = 1514.091
$ws.Range("I122").Value = 1415.5
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4246.5
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1796.5
$ws.Range("N122").Value = -12400

$ws.Range("H129").Value = 2124.5293
$ws.Range("I129").Value = 1863.9166
$ws.Range("J129").Value = 2750
$ws.Range("K129").Value = 5591.7498
$ws.Range("L129").Value = 8250
$ws.Range("M129").Value = -591.7497999999996
$ws.Range("N129").Value = -18250

$ws.Range("H137").Value = 45456076
$ws.Range("I137").Value = 90910000
$ws.Range("J137").Value = 2152
$ws.Range("K137").Value = 272730000
$ws.Range("L137").Value = 6456
$ws.Range("M137").Value = -272727450
$ws.Range("N137").Value = -11556

$ws.Range("H138").Value = 7226.4565
$ws.Range("I138").Value = 11539.682
$ws.Range("J138").Value = 3272.6667
$ws.Range("K138").Value = 34619.046
$ws.Range("L138").Value = 9818.000100000001
$ws.Range("M138").Value = -29479.046
$ws.Range("N138").Value = -20098.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 254
$ws.Range("I5").Value = 305
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 305
$ws.Range("L5").Value = 50
$ws.Range("M5").Value = -193
$ws.Range("N5").Value = -274

$ws.Range("H32").Value = 1487.8572
$ws.Range("I32").Value = 1335.6061
$ws.Range("K32").Value = 1335.6061
$ws.Range("M32").Value = -1048.6061

$ws.Range("H74").Value = 3310.9666
$ws.Range("I74").Value = 1691.9524
$ws.Range("J74").Value = 7088.6665
$ws.Range("K74").Value = 1691.9524
$ws.Range("L74").Value = 7088.6665
$ws.Range("M74").Value = -817.9523999999999
$ws.Range("N74").Value = -8836.666499999999

$ws.Range("H77").Value = 3310.9666
$ws.Range("I77").Value = 1691.9524
$ws.Range("J77").Value = 7088.6665
$ws.Range("K77").Value = 8459.761999999999
$ws.Range("L77").Value = 35443.3325
$ws.Range("M77").Value = -4091.761999999999
$ws.Range("N77").Value = -44179.3325

$ws.Range("H113").Value = 98596.5
$ws.Range("J113").Value = 98596.5
$ws.Range("L113").Value = 98596.5
$ws.Range("N113").Value = -107274.5

$ws.Range("H132").Value = 2716.4546
$ws.Range("I132").Value = 1783.5
$ws.Range("K132").Value = 5350.5
$ws.Range("M132").Value = -2820.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 254
$ws.Range("I4").Value = 305
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 305
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = -190
$ws.Range("N4").Value = -280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 20030000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H132").Value = 54058580
$ws.Range("I132").Value = 83336904
$ws.Range("K132").Value = 250010712
$ws.Range("M132").Value = -250008182

$ws.Range("H134").Value = 2474.4285
$ws.Range("I134").Value = 2085.2354
$ws.Range("J134").Value = 4128.5
$ws.Range("K134").Value = 6255.706200000001
$ws.Range("L134").Value = 12385.5
$ws.Range("M134").Value = -3720.706200000001
$ws.Range("N134").Value = -17455.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 66.44444
$ws.Range("I2").Value = 11.5
$ws.Range("K2").Value = 69
$ws.Range("M2").Value = 44

$ws.Range("H7").Value = 150166.3
$ws.Range("I7").Value = 110.5
$ws.Range("J7").Value = 375250
$ws.Range("K7").Value = 331.5
$ws.Range("L7").Value = 1125750
$ws.Range("M7").Value = -219.5
$ws.Range("N7").Value = -1125974

$ws.Range("H121").Value = 33338020
$ws.Range("I121").Value = 66667040
$ws.Range("J121").Value = 9000
$ws.Range("K121").Value = 200001120
$ws.Range("L121").Value = 27000
$ws.Range("M121").Value = -199999810
$ws.Range("N121").Value = -29620

$ws.Range("H132").Value = 1424.091
$ws.Range("I132").Value = 1027.8572
$ws.Range("K132").Value = 9250.7148
$ws.Range("M132").Value = -6720.7148

$ws.Range("H139").Value = 68085.60000000001
$ws.Range("I139").Value = 72591.71000000001
$ws.Range("K139").Value = 217775.13
$ws.Range("M139").Value = -212635.13

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1003180
$ws.Range("I80").Value = 1252476.2
$ws.Range("K80").Value = 1252476.2
$ws.Range("M80").Value = -1251478.2

$ws.Range("H83").Value = 1003180
$ws.Range("I83").Value = 1252476.2
$ws.Range("K83").Value = 6262381
$ws.Range("M83").Value = -6257389

$ws.Range("H122").Value = 47530
$ws.Range("I122").Value = 75348.5
$ws.Range("K122").Value = 226045.5
$ws.Range("M122").Value = -223595.5

$ws.Range("H132").Value = 11116684
$ws.Range("I132").Value = 16398986
$ws.Range("J132").Value = 5636.3105
$ws.Range("K132").Value = 49196958
$ws.Range("L132").Value = 16908.9315
$ws.Range("M132").Value = -49194428
$ws.Range("N132").Value = -21968.9315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3766
$ws.Range("I40").Value = 3611.8125
$ws.Range("J40").Value = 4999.5
$ws.Range("K40").Value = 3611.8125
$ws.Range("L40").Value = 4999.5
$ws.Range("M40").Value = -3475.8125
$ws.Range("N40").Value = -5271.5

$ws.Range("H82").Value = 2531.8572
$ws.Range("I82").Value = 799.7143
$ws.Range("J82").Value = 3397.9285
$ws.Range("K82").Value = 799.7143
$ws.Range("L82").Value = 3397.9285
$ws.Range("M82").Value = -438.7143
$ws.Range("N82").Value = -4119.9285

$ws.Range("H85").Value = 2531.8572
$ws.Range("I85").Value = 799.7143
$ws.Range("J85").Value = 3397.9285
$ws.Range("K85").Value = 799.7143
$ws.Range("L85").Value = 3397.9285
$ws.Range("M85").Value = 448.2857
$ws.Range("N85").Value = -5893.9285

$ws.Range("H122").Value = 5443.75
$ws.Range("I122").Value = 4827.273
$ws.Range("K122").Value = 14481.819
$ws.Range("M122").Value = -12031.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H107").Value = 1791.3334
$ws.Range("I107").Value = 911
$ws.Range("J107").Value = 4306.5713
$ws.Range("K107").Value = 2733
$ws.Range("L107").Value = 12919.7139
$ws.Range("M107").Value = -813
$ws.Range("N107").Value = -16759.7139

$ws.Range("H122").Value = 2637.5454
$ws.Range("I122").Value = 2446.7632
$ws.Range("J122").Value = 3845.8333
$ws.Range("K122").Value = 7340.2896
$ws.Range("L122").Value = 11537.4999
$ws.Range("M122").Value = -4890.2896
$ws.Range("N122").Value = -16437.4999

$ws.Range("H135").Value = 88018.42999999999
$ws.Range("J135").Value = 88018.42999999999
$ws.Range("L135").Value = 88018.42999999999
$ws.Range("N135").Value = -98158.42999999999

